$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the "NA" text value into F6, F7, F8, F9, F16 (ETA column), preserving existing cell style/format
$ws.Range("F6").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("F16").Value = "NA"

# Update the view: scroll so row 4 is the top-left visible row, and move the
# active selection to F17 (a single cell, just past the last data row)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("F17").Select()
